$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.496.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.914.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4834"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2893"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "109.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.917.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6724"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.511.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007568"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.168.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.497"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.460"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.404"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.153"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04992"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7308"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.726"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.669"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.018"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4455"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8652"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.809"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.352"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.243"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1240"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
